$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value  = "Planifier des rendez-vous avec le client"
$ws.Range("C5").Value  = "Toujours être en contact sur l'évolution et échanger au sein de l'équipe"
$ws.Range("C7").Value  = "Les documents doivent être prêt à temps et complets"
$ws.Range("C9").Value  = "Un code robuste sera plus simple a faire évoluer "
$ws.Range("C10").Value = "Le formats de données devra être choisis par l'équipe"
$ws.Range("C6").Value  = "Chaque membres doit être assidus et répondrent au planning mis en place"
$ws.Range("C12").Value = "Nous devrons utiliser chacunes des compétences technologiques présentent au sein de l'équipe "
$ws.Range("C13").Value = "Le projet devra être bien compris pour rendre quelque chose de complet en adéquation avec la demande client"
$ws.Range("C15").Value = "Connaître les différents point que comporteront les devis"
$ws.Range("C16").Value = "Mise au point sur les différents moyens de sauvegarde"
$ws.Range("C17").Value = "Mise au point avec le client sur le besoin d'être en ligne ou hors ligne"
$ws.Range("C19").Value = "Des réductions devront pouvoir être appliquées aux devis"

$ws.Range("B18").Select()
